# Generate Report for Handoff
# Updates the localization status report: all "In Translation" rows are now
# "Ready for handoff", and the relevant "Latest ... Date(time)" timestamps
# are refreshed to reflect this handoff generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status columns: "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# Timestamps refreshed for this handoff generation run
$overview.Range("G2").Value = "2016-09-05 03:06:27"
$dede.Range("H2").Value     = "2016-09-05 03:06:27"
$zhcn.Range("H2").Value     = "2016-09-05 03:06:23"

# The longer "Ready for handoff" text widens the Status columns accordingly
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333
